$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet and rename it
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $cell = $newSheet.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$data = @(
    @("4167", $null, "", "", "", "NO"),
    @("4168", 1, "1", "1", "6.28%", "NO"),
    @("4170", 7, "4", "0", "22.93%", "NO"),
    @("4222", $null, "", "", "", "NO"),
    @("4415", 7, "1", "0", "8.29%", "NO"),
    @("4419", $null, "", "", "", "NO"),
    @("4421", $null, "", "", "", "NO"),
    @("4423", 4, "0", "0", "1.94%", "NO")
)

# Columns A, C, D, E, F hold text values (even when numeric-looking, or
# empty text when blank); column B holds true numbers when present and
# empty text when blank. A lone "'" forces text typing on an otherwise
# empty value (mirrors the leading-apostrophe "treat as text" convention)
# without leaving the cell completely empty/untyped.
for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $cell = $newSheet.Cells.Item($r + 2, $c + 1)
        $value = $row[$c]
        if ($c -eq 1) {
            if ($null -eq $value) {
                $cell.NumberFormat = "@"
                $cell.Value = "'"
                $cell.Style = "Normal"
            } else {
                $cell.Value = $value
            }
        } else {
            $cell.NumberFormat = "@"
            if ($value -eq "") {
                $cell.Value = "'"
            } else {
                $cell.Value = $value
            }
            $cell.Style = "Normal"
        }
    }
}
